# Update stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = 1002
$ws.Range("D24").Value = 5966728
$ws.Range("E24").Value = 934.4914643696163
$ws.Range("G24").Value = 3.83419689119171
$ws.Range("H24").Value = 26.39841104001803
